$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.614.52'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.923.40'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''246.63'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '''0.4741'
$ws.Range('D8').Value = '''0.2883'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D10').Value = '''105.20'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '''18.33'
$ws.Range('E11').Value = '  -4.21%  '
$ws.Range('D12').Value = '1.920.96'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '''0.07694'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = '''5.337'
$ws.Range('E14').Value = '  +4.12%  '
$ws.Range('D15').Value = '''0.6678'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').Value = '''291.65'
$ws.Range('D17').Value = '30.612.80'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '''0.000007617'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').Value = '''0.9995'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''12.96'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''5.552'
$ws.Range('E21').Value = '  +5.94%  '
$ws.Range('D22').Value = '2.170.60'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '''6.427'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('D25').Value = '''9.460'
$ws.Range('E25').Value = '  +2.90%  '
$ws.Range('D26').Value = '''167.52'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +6.50%  '
$ws.Range('D28').Value = '''2.113'
$ws.Range('E28').Value = '  +5.53%  '
$ws.Range('D29').Value = '''0.1072'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('D30').Value = '''1.398'
$ws.Range('E30').Value = '  +3.61%  '
$ws.Range('D31').Value = '''4.179'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('D32').Value = '''4.059'
$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('D33').Value = '''0.05037'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').Value = '''0.7381'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').Value = '''1.143'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '''0.02063'
$ws.Range('E36').Value = '  +6.13%  '
$ws.Range('D37').Value = '''2.740'
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').Value = '''2.687'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').Value = '''111.28'
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('D41').Value = '''0.8738'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '''0.4365'
$ws.Range('E42').Value = '  +5.91%  '
$ws.Range('D43').Value = '''5.915'
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '''67.98'
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').Value = '''7.271'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '''9.297'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('D48').Value = '''48.15'
$ws.Range('E48').Value = '  +14.61%  '
$ws.Range('D49').Value = '''0.1240'
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''34.99'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.2509'
$ws.Range('E51').Value = '  +12.10%  '
